# Update the "Förändrad" (changed) date column C for rows 2-45
# from serial date 45179 (2023-09-10) to 45180 (2023-09-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2:C45").Value = 45180
